# The section header row "grandes regiões e unidades da federação" (row 6)
# was a stray/duplicate label with no data beneath it. The correction removes
# that empty header row entirely so the "norte" data row (previously row 7)
# and everything below it shifts up by one, and the trailing empty slot
# (previously row 38) disappears.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
